# This script applies the corrections to the SD1 sheet:
# - Column A (Numero): change academic year from 2015 to 2017 for all data rows (3-63)
# - Column E (Moyenne de l'etudiant): update grade values per corrected dataset
#
# Context (commit message): correction of JSON data generation when a student's
# data (UE, semester, or training year) does not exist in the database -> replaced
# by NULL. The year values in column A and some averages in column E were fixed
# accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("A3").Value = 20170926
$ws.Range("E3").Value = 13
$ws.Range("A4").Value = 20170927
$ws.Range("E4").Value = 10
$ws.Range("A5").Value = 20170928
$ws.Range("A6").Value = 20170929
$ws.Range("E6").Value = 11
$ws.Range("A7").Value = 20170930
$ws.Range("E7").Value = 17
$ws.Range("A8").Value = 20170931
$ws.Range("E8").Value = 19
$ws.Range("A9").Value = 20170932
$ws.Range("E9").Value = 11
$ws.Range("A10").Value = 20170933
$ws.Range("E10").Value = 15
$ws.Range("A11").Value = 20170934
$ws.Range("E11").Value = 18
$ws.Range("A12").Value = 20170935
$ws.Range("E12").Value = 6
$ws.Range("A13").Value = 20170936
$ws.Range("E13").Value = 18
$ws.Range("A14").Value = 20170937
$ws.Range("E14").Value = 17
$ws.Range("A15").Value = 20170938
$ws.Range("E15").Value = 13
$ws.Range("A16").Value = 20170939
$ws.Range("E16").Value = 14
$ws.Range("A17").Value = 20170940
$ws.Range("E17").Value = 20
$ws.Range("A18").Value = 20170941
$ws.Range("E18").Value = 6
$ws.Range("A19").Value = 20170942
$ws.Range("E19").Value = 9
$ws.Range("A20").Value = 20170943
$ws.Range("E20").Value = 18
$ws.Range("A21").Value = 20170944
$ws.Range("E21").Value = 12
$ws.Range("A22").Value = 20170945
$ws.Range("E22").Value = 7
$ws.Range("A23").Value = 20170946
$ws.Range("E23").Value = 19
$ws.Range("A24").Value = 20170947
$ws.Range("E24").Value = 14
$ws.Range("A25").Value = 20170948
$ws.Range("E25").Value = 12
$ws.Range("A26").Value = 20170949
$ws.Range("E26").Value = 13
$ws.Range("A27").Value = 20170950
$ws.Range("E27").Value = 16
$ws.Range("A28").Value = 20170951
$ws.Range("A29").Value = 20170952
$ws.Range("E29").Value = 5
$ws.Range("A30").Value = 20170953
$ws.Range("E30").Value = 14
$ws.Range("A31").Value = 20170954
$ws.Range("E31").Value = 12
$ws.Range("A32").Value = 20170955
$ws.Range("E32").Value = 5
$ws.Range("A33").Value = 20170956
$ws.Range("A34").Value = 20170957
$ws.Range("E34").Value = 9
$ws.Range("A35").Value = 20170958
$ws.Range("A36").Value = 20170959
$ws.Range("E36").Value = 6
$ws.Range("A37").Value = 20170960
$ws.Range("E37").Value = 11
$ws.Range("A38").Value = 20170961
$ws.Range("E38").Value = 9
$ws.Range("A39").Value = 20170962
$ws.Range("E39").Value = 14
$ws.Range("A40").Value = 20170963
$ws.Range("E40").Value = 13
$ws.Range("A41").Value = 20170964
$ws.Range("E41").Value = 12
$ws.Range("A42").Value = 20170965
$ws.Range("E42").Value = 11
$ws.Range("A43").Value = 20170966
$ws.Range("E43").Value = 19
$ws.Range("A44").Value = 20170967
$ws.Range("E44").Value = 12
$ws.Range("A45").Value = 20170968
$ws.Range("E45").Value = 8
$ws.Range("A46").Value = 20170969
$ws.Range("E46").Value = 9
$ws.Range("A47").Value = 20170970
$ws.Range("E47").Value = 5
$ws.Range("A48").Value = 20170971
$ws.Range("E48").Value = 14
$ws.Range("A49").Value = 20170972
$ws.Range("E49").Value = 18
$ws.Range("A50").Value = 20170973
$ws.Range("E50").Value = 14
$ws.Range("A51").Value = 20170974
$ws.Range("E51").Value = 14
$ws.Range("A52").Value = 20170975
$ws.Range("E52").Value = 6
$ws.Range("A53").Value = 20170976
$ws.Range("E53").Value = 8
$ws.Range("A54").Value = 20170977
$ws.Range("E54").Value = 16
$ws.Range("A55").Value = 20170978
$ws.Range("E55").Value = 8
$ws.Range("A56").Value = 20170979
$ws.Range("E56").Value = 13
$ws.Range("A57").Value = 20170980
$ws.Range("E57").Value = 10
$ws.Range("A58").Value = 20170981
$ws.Range("E58").Value = 19
$ws.Range("A59").Value = 20170982
$ws.Range("E59").Value = 5
$ws.Range("A60").Value = 20170983
$ws.Range("E60").Value = 7
$ws.Range("A61").Value = 20170984
$ws.Range("E61").Value = 16
$ws.Range("A62").Value = 20170985
$ws.Range("E62").Value = 17
$ws.Range("A63").Value = 20170986
$ws.Range("E63").Value = 12
